# unify the conception of DataNode, DataTable, Entity.
# Rename the worksheet from "Property1" to "DataNode" and move the
# active selection to D42 (matching the resaved workbook's cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "DataNode"

$ws.Range("D42").Select()
